# Edit script generated to apply the documented diff.
$doc = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1) Date field: "28/09/2020" -> "09" + [_GoBack bookmark] + "/11" + "/2020"
# ----------------------------------------------------------------------

# The _GoBack bookmark currently sits elsewhere in the document (inside the
# "Un programa en C..." paragraph). Remove it there; it will be re-created
# at its new location (inside the date field) below.
if ($doc.Bookmarks.Exists("_GoBack")) {
    $doc.Bookmarks("_GoBack").Delete()
}

# Replace the date text.
$dateRng = $doc.Content
$dateRng.Find.Execute("28/09/2020", $false, $false, $false, $false, $false, $true, 1, $false, "09/11/2020", 2)

# Locate the new date text so we can split it into three runs and drop a
# bookmark between the first and second run.
$dateRng2 = $doc.Content
$dateRng2.Find.Execute("09/11/2020", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dateStart = $dateRng2.Start

# Force a run boundary between "09" and "/11" and between "/11" and "/2020"
# by toggling Bold off/on (a no-op visually) which causes the host to split
# the run at those character offsets.
$splitRng = $doc.Range($dateStart + 2, $dateStart + 5)
$splitRng.Bold = 1
$splitRng.Bold = 0

# Insert the (now-empty, i.e. collapsed) bookmark between "09" and "/11".
$bmRng = $doc.Range($dateStart + 2, $dateStart + 2)
$doc.Bookmarks.Add("_GoBack", $bmRng)

# ----------------------------------------------------------------------
# 2) Merge the runs in the "Éste es muy utilizado..." paragraph into one run.
# ----------------------------------------------------------------------
$run1Rng1 = $doc.Content
$run1Rng1.Find.Execute("Éste es muy utilizado ya que la forma de dar instrucciones es muy cercana a lo que un humano podría abstraer, es decir, las instrucciones no son tal cual las que una computadora podría entender", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$run1Start1 = $run1Rng1.Start
$run1End1 = $run1Rng1.End

$fullRng1 = $doc.Content
$fullRng1.Find.Execute("Éste es muy utilizado ya que la forma de dar instrucciones es muy cercana a lo que un humano podría abstraer, es decir, las instrucciones no son tal cual las que una computadora podría entender. Por esta razón, C es conocido como un lenguaje de alto nivel, esto significa a que las instrucciones podrían ser entendidas fácilmente por un humano. Algunos autores consideran al lenguaje C como un lenguaje de mediano nivel, ya que no es totalmente transparente sino tiene elementos que tienen que ver con la arquitectura de la máquina a la hora de programar.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fullEnd1 = $fullRng1.End

# Delete the trailing runs (everything after run 1's text, up to the end of
# the merged text) then retype the remainder onto the back of run 1 so it
# absorbs the text and keeps run 1's original formatting/identity.
$trailRng1 = $doc.Range($run1End1, $fullEnd1)
$trailRng1.Delete()

$run1RngAgain1 = $doc.Range($run1Start1, $run1End1)
$run1RngAgain1.InsertAfter(". Por esta razón, C es conocido como un lenguaje de alto nivel, esto significa a que las instrucciones podrían ser entendidas fácilmente por un humano. Algunos autores consideran al lenguaje C como un lenguaje de mediano nivel, ya que no es totalmente transparente sino tiene elementos que tienen que ver con la arquitectura de la máquina a la hora de programar.")

# ----------------------------------------------------------------------
# 3) Merge the runs in the "Un programa en C..." paragraph into one run
#    (this also removes the _GoBack bookmark that used to live here).
# ----------------------------------------------------------------------
$run1Rng2 = $doc.Content
$run1Rng2.Find.Execute("Un programa en C se elabora describ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$run1Start2 = $run1Rng2.Start
$run1End2 = $run1Rng2.End

$fullRng2 = $doc.Content
$fullRng2.Find.Execute("Un programa en C se elabora describiendo cada una de las instrucciones de acuerdo a las reglas definidas en este lenguaje en un archivo de texto para después ser procesadas en un compilador. Un compilador es un programa que toma como entrada un archivo de texto y tiene como salida un programa ejecutable, éste tiene instrucciones que poden ser procesadas por el hardware de la computadora en conjunto con el sistema operativo que corre sobre ella. Se tiene como ventaja que un programa escrito en lenguaje C, siguiendo siempre su estándar, puede correr en cualquier máquina siempre y cuando exista un compilador de C hecho para tal.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fullEnd2 = $fullRng2.End

# Delete the trailing runs (everything after run 1's text, up to the end of
# the merged text) then retype the remainder onto the back of run 1 so it
# absorbs the text and keeps run 1's original formatting/identity.
$trailRng2 = $doc.Range($run1End2, $fullEnd2)
$trailRng2.Delete()

$run1RngAgain2 = $doc.Range($run1Start2, $run1End2)
$run1RngAgain2.InsertAfter("iendo cada una de las instrucciones de acuerdo a las reglas definidas en este lenguaje en un archivo de texto para después ser procesadas en un compilador. Un compilador es un programa que toma como entrada un archivo de texto y tiene como salida un programa ejecutable, éste tiene instrucciones que poden ser procesadas por el hardware de la computadora en conjunto con el sistema operativo que corre sobre ella. Se tiene como ventaja que un programa escrito en lenguaje C, siguiendo siempre su estándar, puede correr en cualquier máquina siempre y cuando exista un compilador de C hecho para tal.")

# ----------------------------------------------------------------------
# 4) Merge the runs in the "Cuando el compilador señala..." paragraph into one run.
# ----------------------------------------------------------------------
$run1Rng3 = $doc.Content
$run1Rng3.Find.Execute("Cuando el compilador señala un error no cabe más que invocar algún editor de texto, revisar cuidadosamente el programa y corregir. Se debe verificar la coherencia total del programa para evitar tener que volver a repetir este paso de manera continua.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$run1Start3 = $run1Rng3.Start
$run1End3 = $run1Rng3.End

$fullRng3 = $doc.Content
$fullRng3.Find.Execute("Cuando el compilador señala un error no cabe más que invocar algún editor de texto, revisar cuidadosamente el programa y corregir. Se debe verificar la coherencia total del programa para evitar tener que volver a repetir este paso de manera continua. A veces el compilador arroja advertencias durante el proceso, se generará el archivo ejecutable, pero puede tener problemas a la hora de ejecución por lo que es mejor investigar de qué tratan o porqué se generaron.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fullEnd3 = $fullRng3.End

# Delete the trailing runs (everything after run 1's text, up to the end of
# the merged text) then retype the remainder onto the back of run 1 so it
# absorbs the text and keeps run 1's original formatting/identity.
$trailRng3 = $doc.Range($run1End3, $fullEnd3)
$trailRng3.Delete()

$run1RngAgain3 = $doc.Range($run1Start3, $run1End3)
$run1RngAgain3.InsertAfter(" A veces el compilador arroja advertencias durante el proceso, se generará el archivo ejecutable, pero puede tener problemas a la hora de ejecución por lo que es mejor investigar de qué tratan o porqué se generaron.")

# ----------------------------------------------------------------------
# 5) Merge the runs in the "La ejecución es la etapa..." paragraph into one run.
# ----------------------------------------------------------------------
$run1Rng4 = $doc.Content
$run1Rng4.Find.Execute("La ejecución es la etapa que sigue después de haber compilado el programa. Uva vez compilado el programa, se puede distribuir para equipos que ejecuten el mismo sistema operativo y tengan la misma plataforma de hardware (tipo de procesador, set de instrucciones y arquitectura en general). Los pasos para realizar la ejecución dependen del sistema operativo y del entorno.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$run1Start4 = $run1Rng4.Start
$run1End4 = $run1Rng4.End

$fullRng4 = $doc.Content
$fullRng4.Find.Execute("La ejecución es la etapa que sigue después de haber compilado el programa. Uva vez compilado el programa, se puede distribuir para equipos que ejecuten el mismo sistema operativo y tengan la misma plataforma de hardware (tipo de procesador, set de instrucciones y arquitectura en general). Los pasos para realizar la ejecución dependen del sistema operativo y del entorno. Es mejor ejecutar el programa en el símbolo de sistema porque, aunque el programa finalice su ejecución, los resultados continuarán siendo visibles en la consola.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fullEnd4 = $fullRng4.End

# Delete the trailing runs (everything after run 1's text, up to the end of
# the merged text) then retype the remainder onto the back of run 1 so it
# absorbs the text and keeps run 1's original formatting/identity.
$trailRng4 = $doc.Range($run1End4, $fullEnd4)
$trailRng4.Delete()

$run1RngAgain4 = $doc.Range($run1Start4, $run1End4)
$run1RngAgain4.InsertAfter(" Es mejor ejecutar el programa en el símbolo de sistema porque, aunque el programa finalice su ejecución, los resultados continuarán siendo visibles en la consola.")

